$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric need to be forced to Text format first,
# so Excel keeps them as the literal strings from the source data feed
# (matching the original inlineStr cells) instead of coercing to numbers.
$textCells = @("D5", "D7", "D11", "D13", "D14", "D19", "D20", "D22", "D25", "D26", "D28", "D30", "D38", "D40", "D41", "D42", "D43", "D47")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '37.837.93'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '2.084.44'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '233.60'
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("D7").Value = '59.32'
$ws.Range("E7").Value = '  +2.75%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("E10").Value = '  +0.61%  '
$ws.Range("D11").Value = '0.107'
$ws.Range("E11").Value = '  +1.48%  '
$ws.Range("D12").Value = '2.392.04'
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").Value = '14.76'
$ws.Range("E13").Value = '  +2.08%  '
$ws.Range("D14").Value = '21.26'
$ws.Range("E14").Value = '  +0.02%  '
$ws.Range("E15").Value = '  +1.06%  '
$ws.Range("E16").Value = '  +1.50%  '
$ws.Range("D17").Value = '2.085.53'
$ws.Range("E17").Value = '  -0.55%  '
$ws.Range("D18").Value = '37.781.12'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").Value = '6.15'
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("D20").Value = '71.60'
$ws.Range("E20").Value = '  +0.92%  '
$ws.Range("D21").Value = '0.0₃0850'
$ws.Range("E21").Value = '  +3.26%  '
$ws.Range("D22").Value = '227.95'
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("E24").Value = '  -0.82%  '
$ws.Range("D25").Value = '2.40'
$ws.Range("E25").Value = '  +0.38%  '
$ws.Range("D26").Value = '171.75'
$ws.Range("E26").Value = '  +0.60%  '
$ws.Range("E27").Value = '  +3.48%  '
$ws.Range("D28").Value = '0.137'
$ws.Range("E28").Value = '  -2.38%  '
$ws.Range("E29").Value = '  -1.44%  '
$ws.Range("D30").Value = '19.50'
$ws.Range("E31").Value = '  +1.50%  '
$ws.Range("E32").Value = '  +2.12%  '
$ws.Range("E33").Value = '  +0.42%  '
$ws.Range("E34").Value = '  +2.32%  '
$ws.Range("E35").Value = '  -0.84%  '
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("E37").Value = '  -0.58%  '
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("E39").Value = '  -0.77%  '
$ws.Range("D40").Value = '0.0986'
$ws.Range("E40").Value = '  -1.65%  '
$ws.Range("D41").Value = '99.48'
$ws.Range("E41").Value = '  +2.03%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.0219'
$ws.Range("E42").Value = '  +2.45%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = '17.13'
$ws.Range("E43").Value = '  +8.37%  '
$ws.Range("E44").Value = '  -1.38%  '
$ws.Range("D45").Value = '1.450.91'
$ws.Range("E45").Value = '  -0.31%  '
$ws.Range("E46").Value = '  -1.16%  '
$ws.Range("D47").Value = '4.16'
$ws.Range("E47").Value = '  +3.36%  '
$ws.Range("E48").Value = '  +0.37%  '
$ws.Range("E49").Value = '  -0.42%  '
$ws.Range("E50").Value = '  -1.01%  '
$ws.Range("D51").Value = '2.277.35'
$ws.Range("E51").Value = '  -0.38%  '
